$d = $word.ActiveDocument

# Replace the first bullet's text with the new first task
$d.Content.Find.Execute("Spawn asteroids, move and rotate.", $false, $false, $false, $false, $false, $true, 1, $false, "Create a quick title screen.", 2)

# Insert the remaining new to-do items as list paragraphs right after the
# first bullet, reusing its list/style formatting via InsertParagraphAfter.
$newItems = @(
    "Create a Game Scene",
    "Add player collision so that player can die.",
    "Add asteroid collision so that player can collide with player.",
    "Have bullets destroy asteroids.",
    "Have asteroid break into to when hit with laser."
)

$idx = 3
foreach ($item in $newItems) {
    $p = $d.Paragraphs($idx)
    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    $newP = $d.Paragraphs($idx)
    $newP.Range.InsertBefore($item)
}

# Replace the URL at the bottom of the document
$d.Content.Find.Execute("http://www.freeasteroids.org/", $false, $false, $false, $false, $false, $true, 1, $false, "http://my.ign.com/atari/asteroids", 2)
